$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HTLP")

# --- Block 1 (rows 3-9): new "Fall" (column E) raw observations ---
$ws.Range("E3").Value = 440.16
$ws.Range("E4").Value = 263.79000000000002

$ws.Range("D7").Value = 589.66999999999996
$ws.Range("E7").Value = 494.82

$ws.Range("D8").Value = 205.92
$ws.Range("E8").Value = 232.18

# --- Block 2 (rows 15-21): new "Fall" (column E) raw observations ---
$ws.Range("E15").Value = 1821.15
$ws.Range("E16").Value = 32.659999999999997

$ws.Range("D19").Value = 1597.37
$ws.Range("E19").Value = 1796.33

$ws.Range("D20").Value = -19.940000000000001
$ws.Range("E20").Value = -11.33

# Move the active cell/selection to F7, matching the saved view state
$ws.Activate()
$ws.Range("F7").Select()
